# Automatic BRVM update (GitHub Actions) - applies refreshed
# "Recommandations" and "Top_YTD" figures to the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Recommandations"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2 - BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Range("D2").Value = 656.8200000000001
$ws1.Range("E2").Value = 163.99

# Row 3 - BRVM - SERVICES FINANCIERS
$ws1.Range("D3").Value = 585.66
$ws1.Range("E3").Value = 147.99

# Row 4 - BRVM-PRESTIGE
$ws1.Range("D4").Value = 570.45
$ws1.Range("E4").Value = 142.94

# Row 5 - BRVM - INDUSTRIELS
$ws1.Range("D5").Value = 549.91
$ws1.Range("E5").Value = 140.11

# Row 6 - BRVM - ENERGIE
$ws1.Range("D6").Value = 451.41
$ws1.Range("E6").Value = 113.51

# Row 7 - BRVM - SERVICES PUBLICS
$ws1.Range("D7").Value = 435.14
$ws1.Range("E7").Value = 110.41

# Row 8 - BRVM - TELECOMMUNICATIONS
$ws1.Range("D8").Value = 372.26
$ws1.Range("E8").Value = 92.78

# Row 12 - EVIOSYS PACKAGING SIEM CI (SEMC)
$ws1.Range("B12").Value = 3
$ws1.Range("C12").Value = 1
$ws1.Range("D12").Value = 18.1
$ws1.Range("E12").Value = -3.08

# Row 13 - now SICABLE CI (CABC)
$ws1.Range("A13").Value = "SICABLE CI (CABC)"
$ws1.Range("D13").Value = 11.18
$ws1.Range("E13").Value = 7.36

# Row 14 - NEI-CEDA CI (NEIC)
$ws1.Range("C14").Value = 1
$ws1.Range("D14").Value = 9.77
$ws1.Range("E14").Value = -4.17
$ws1.Range("G14").Value = "👀 À surveiller"

# Row 15 - now UNILEVER CI (UNLC)
$ws1.Range("A15").Value = "UNILEVER CI (UNLC)"
$ws1.Range("D15").Value = 8.33
$ws1.Range("E15").Value = 7.5

# Row 16 - now FILTISAC CI (FTSC)
$ws1.Range("A16").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B16").Value = 1
$ws1.Range("C16").Value = 0
$ws1.Range("D16").Value = 7.26
$ws1.Range("E16").Value = 7.26
$ws1.Range("G16").Value = "➖ Neutre"

# Row 17 - now BANK OF AFRICA BN (BOAB)
$ws1.Range("A17").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("B17").Value = 1
$ws1.Range("C17").Value = 0
$ws1.Range("D17").Value = 6.19
$ws1.Range("E17").Value = 6.19
$ws1.Range("G17").Value = "➖ Neutre"

# Row 18 - now SOLIBRA CI (SLBC)
$ws1.Range("A18").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B18").Value = 2
$ws1.Range("C18").Value = 2
$ws1.Range("D18").Value = 5.29
$ws1.Range("E18").Value = -1.99
$ws1.Range("G18").Value = "👀 À surveiller"

# Row 19 - now CFAO MOTORS CI (CFAC)
$ws1.Range("A19").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("C19").Value = 0
$ws1.Range("D19").Value = 4.71
$ws1.Range("E19").Value = 4.71
$ws1.Range("G19").Value = "➖ Neutre"

# Row 20 - now SICOR CI (SICC)
$ws1.Range("A20").Value = "SICOR CI (SICC)"
$ws1.Range("B20").Value = 2
$ws1.Range("D20").Value = 4.24
$ws1.Range("E20").Value = 4.08

# Row 21 - now ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Range("A21").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("D21").Value = 0.2
$ws1.Range("E21").Value = 4.55

# Row 22 - now SETAO CI (STAC)
$ws1.Range("A22").Value = "SETAO CI (STAC)"
$ws1.Range("B22").Value = 2
$ws1.Range("C22").Value = 2
$ws1.Range("D22").Value = 0.02
$ws1.Range("E22").Value = 5.42
$ws1.Range("G22").Value = "👀 À surveiller"

# Row 23 - now ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)
$ws1.Range("A23").Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Range("B23").Value = 1
$ws1.Range("D23").Value = -0.15
$ws1.Range("E23").Value = 3.75
$ws1.Range("G23").Value = "👀 À surveiller"

# Row 24 - now VIVO ENERGY CI (SHEC)
$ws1.Range("A24").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("D24").Value = -1.29
$ws1.Range("E24").Value = -1.29

# Row 25 - now ONATEL BF (ONTBF)
$ws1.Range("A25").Value = "ONATEL BF (ONTBF)"
$ws1.Range("B25").Value = 0
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = -1.43
$ws1.Range("E25").Value = -1.43
$ws1.Range("G25").Value = "➖ Neutre"

# Row 26 - now TOTALENERGIES MARKETING CI (TTLC)
$ws1.Range("A26").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("D26").Value = -1.49
$ws1.Range("E26").Value = -1.49

# Row 27 - now SUCRIVOIRE (SCRC)
$ws1.Range("A27").Value = "SUCRIVOIRE (SCRC)"
$ws1.Range("D27").Value = -1.9
$ws1.Range("E27").Value = -1.9

# Row 28 - now NESTLE CI (NTLC)
$ws1.Range("A28").Value = "NESTLE CI (NTLC)"
$ws1.Range("D28").Value = -2.55
$ws1.Range("E28").Value = -2.55

# Row 29 - now SOCIETE IVOIRIENNE DE BANQUE  (SIBC)
$ws1.Range("A29").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("D29").Value = -2.68
$ws1.Range("E29").Value = -2.68

# Row 30 - now SAFCA CI (SAFC)
$ws1.Range("A30").Value = "SAFCA CI (SAFC)"
$ws1.Range("D30").Value = -2.69
$ws1.Range("E30").Value = -2.69

# Row 31 - now SAPH CI (SPHC)
$ws1.Range("A31").Value = "SAPH CI (SPHC)"
$ws1.Range("B31").Value = 0
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = -3.56
$ws1.Range("E31").Value = -3.56
$ws1.Range("G31").Value = "➖ Neutre"

# ---------------------------------------------------------------------
# Sheet "Top_YTD"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")

$ws2.Range("B2").Value = 4772.62
$ws2.Range("B3").Value = 3586.81
$ws2.Range("B4").Value = 3364.58
$ws2.Range("B5").Value = 3080.02
$ws2.Range("B6").Value = 1952.63
$ws2.Range("B7").Value = 1800
$ws2.Range("B8").Value = 1289.35
